$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$gval = [double]"-1.20892581961463e+24"

# Update row 2
$ws.Range("B2").Value = -5
$ws.Range("C2").Value = -5
$ws.Range("D2").Value = 40
$ws.Range("E2").Value = 14.9990234375
$ws.Range("F2").Value = 14.9990234375
$ws.Range("G2").Value = $gval
$ws.Range("H2").Value = 1.000005

# Update row 3
$ws.Range("B3").Value = -5
$ws.Range("C3").Value = -5
$ws.Range("D3").Value = 40
$ws.Range("E3").Value = 14.9990234375
$ws.Range("F3").Value = 14.9990234375
$ws.Range("G3").Value = $gval
$ws.Range("H3").Value = 0

# Delete rows 4 through 11 (rows below row 3)
$ws.Range("A4:H11").Delete()
